$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = 1021
$ws.Range("J4").Value = 1106
$ws.Range("Q4").Value = 742
